# Update countries & provincias Spain
# - Refresh the COVID numeric counters (columns B..H) for the countries whose
#   rank changed between the 21:00 and 22:17 pulls.
# - Where two countries swapped rank, the row keeps its position in the table
#   but the country name (column A) and its numbers move together, so we
#   overwrite column A with the new country name for that row as well.
# - Bump the "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 16 de Septiembre de 2020 a las 22:17"

# --- Plain numeric refreshes (no rank change) --------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 6815915
$ws.Range("C4").Value = 27768
$ws.Range("D4").Value = 4095623
$ws.Range("E4").Value = 2519279
$ws.Range("G4").Value = 816
$ws.Range("H4").Value = 201013

# Row 5: India
$ws.Range("B5").Value = 5115893
$ws.Range("C5").Value = 97859
$ws.Range("D5").Value = 4022049
$ws.Range("E5").Value = 1010614
$ws.Range("G5").Value = 1139
$ws.Range("H5").Value = 83230

# Row 11: Sudafrica
$ws.Range("B11").Value = 653444
$ws.Range("C11").Value = 1923
$ws.Range("D11").Value = 584195
$ws.Range("E11").Value = 53544
$ws.Range("G11").Value = 64
$ws.Range("H11").Value = 15705

# Row 25: Alemania
$ws.Range("B25").Value = 266865
$ws.Range("C25").Value = 2021
$ws.Range("E25").Value = 18316
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 9449

# Row 29: Canada
$ws.Range("B29").Value = 139454
$ws.Range("C29").Value = 651
$ws.Range("D29").Value = 122008
$ws.Range("E29").Value = 8253

# --- Nepal / Costa Rica swap (rows 55-56, shared-string idx 59/60) -----
$ws.Range("A55").Value = "Costa Rica"
$ws.Range("B55").Value = 59516
$ws.Range("C55").Value = 1379
$ws.Range("D55").Value = 21752
$ws.Range("E55").Value = 37115
$ws.Range("G55").Value = 16
$ws.Range("H55").Value = 649

$ws.Range("A56").Value = "Nepal"
$ws.Range("B56").Value = 58327
$ws.Range("C56").Value = 1539
$ws.Range("D56").Value = 41706
$ws.Range("E56").Value = 16242
$ws.Range("G56").Value = 8
$ws.Range("H56").Value = 379

# Row 112: Malaui
$ws.Range("B112").Value = 5704
$ws.Range("C112").Value = 3
$ws.Range("D112").Value = 3764
$ws.Range("E112").Value = 1762

# --- Surinam / Ruanda swap (rows 123-124, shared-string idx 127/128) ---
$ws.Range("A123").Value = "Ruanda"
$ws.Range("B123").Value = 4634
$ws.Range("C123").Value = 10
$ws.Range("D123").Value = 2789
$ws.Range("E123").Value = 1823
$ws.Range("H123").Value = 22

$ws.Range("A124").Value = "Surinam"
$ws.Range("B124").Value = 4625
$ws.Range("D124").Value = 3996
$ws.Range("E124").Value = 534
$ws.Range("H124").Value = 95

# --- Trinidad yTobago / Sri Lanka / Aruba 3-way rotation (rows 137-139) -
# Original order: Trinidad yTobago(137), Sri Lanka(138), Aruba(139)
# New order:      Aruba(137), Trinidad yTobago(138), Sri Lanka(139)
$ws.Range("A137").Value = "Aruba"
$ws.Range("B137").Value = 3328
$ws.Range("C137").Value = 176
$ws.Range("D137").Value = 1676
$ws.Range("E137").Value = 1630
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 22

$ws.Range("A138").Value = "Trinidad yTobago"
$ws.Range("B138").Value = 3293
$ws.Range("C138").Value = 70
$ws.Range("D138").Value = 810
$ws.Range("E138").Value = 2426
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 57

$ws.Range("A139").Value = "Sri Lanka"
$ws.Range("B139").Value = 3271
$ws.Range("D139").Value = 3021
$ws.Range("E139").Value = 237
$ws.Range("H139").Value = 13

# Row 151: Sierra Leona
$ws.Range("B151").Value = 2133
$ws.Range("C151").Value = 7
$ws.Range("D151").Value = 1646
$ws.Range("E151").Value = 415

# Row 161: (rank unchanged, no country swap)
$ws.Range("B161").Value = 1483
$ws.Range("C161").Value = 45
$ws.Range("D161").Value = 1054
$ws.Range("E161").Value = 376

# --- Lesoto / Liberia swap (rows 162-163, shared-string idx 166/167) ---
$ws.Range("A162").Value = "Liberia"
$ws.Range("B162").Value = 1332
$ws.Range("C162").Value = 5
$ws.Range("D162").Value = 1214
$ws.Range("E162").Value = 36
$ws.Range("H162").Value = 82

$ws.Range("A163").Value = "Lesoto"
$ws.Range("D163").Value = 687
$ws.Range("E163").Value = 607
$ws.Range("H163").Value = 33

# --- Montserrat / Islas Malvinas swap (rows 214-215, idx 218/219) ------
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
